$d = $word.ActiveDocument
$section = $d.Sections.First
$header = $section.Headers.Item(1)
$header.LinkToPrevious = $false
$header.Range.Text = "Questionnaire 37"
$header.Range.ParagraphFormat.Alignment = 1
$header.Range.Font.Name = "Arial"
$header.Range.Font.Size = 12
